$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 94-106 ---
$ws = $wb.Worksheets.Item("PIR")
$startRow = 94
$endRow = 106
$rng = $ws.Range("A" + $startRow + ":F" + $endRow)
$rng.NumberFormat = "@"

$data = @(
    ,("2026-01-28","15:36:09","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:10","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:15","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:20","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:25","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:30","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:35","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:40","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:45","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:50","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:36:55","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:37:01","15:00","Bathroom","No Motion","Inactive")
    ,("2026-01-28","15:37:06","15:00","Bathroom","No Motion","Inactive")
)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Humidity sheet: append rows 95-108 ---
$ws = $wb.Worksheets.Item("Humidity")
$startRow = 95
$endRow = 108
$rng = $ws.Range("A" + $startRow + ":F" + $endRow)
$rng.NumberFormat = "@"

$data = @(
    ,("2026-01-28","15:36:08","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:09","15:00","Bathroom","87.5%","Active")
    ,("2026-01-28","15:36:11","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:19","15:00","Bathroom","88.3%","Active")
    ,("2026-01-28","15:36:23","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:27","15:00","Bathroom","87.4%","Active")
    ,("2026-01-28","15:36:31","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:39","15:00","Bathroom","87.4%","Active")
    ,("2026-01-28","15:36:43","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:47","15:00","Bathroom","87.5%","Active")
    ,("2026-01-28","15:36:51","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:36:59","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:37:03","15:00","Bathroom","88.4%","Active")
    ,("2026-01-28","15:37:07","15:00","Bathroom","87.5%","Active")
)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Temperature sheet: append rows 95-108 ---
$ws = $wb.Worksheets.Item("Temperature")
$startRow = 95
$endRow = 108
$rng = $ws.Range("A" + $startRow + ":F" + $endRow)
$rng.NumberFormat = "@"

$data = @(
    ,("2026-01-28","15:36:09","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:10","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:11","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:19","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:23","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:27","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:31","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:39","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:43","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:47","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:51","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:36:59","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:37:03","15:00","Bathroom","22.9C","Active")
    ,("2026-01-28","15:37:07","15:00","Bathroom","22.9C","Active")
)

$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

